{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Source paragraph:      \"F\u00c9RIAS\"  (plain, no formatting)\n// Target paragraph:      \"Check list guia para entrevista\"\n//   - paragraph centered\n//   - bold + italic, dark-grey themed color, size 14pt (sz=28 half-points)\n//   - WordArt-style text effects (shadow + outline) carried as w14 run\n//     extensions that have no first-class Office.js/Word-OM property, so\n//     the run/paragraph XML is written directly via insertOoxml().\n//   - a `_GoBack` bookmark sits between \"C\" and \"heck\" (an artifact of\n//     Word leaving the last-edit-position bookmark behind) and the two\n//     words \"Check\"/\"list\" are wrapped in proofErr spell-check markers.\n\n// The run formatting (identical on every run and on the paragraph mark)\n// reproduced from the target OOXML.\nconst rPrInner =\n  '<w:b/><w:i/><w:iCs/>' +\n  '<w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/>' +\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>' +\n  '<w14:shadow w14:blurRad=\"0\" w14:dist=\"38100\" w14:dir=\"2700000\" ' +\n  'w14:sx=\"100000\" w14:sy=\"100000\" w14:kx=\"0\" w14:ky=\"0\" w14:algn=\"bl\">' +\n  '<w14:schemeClr w14:val=\"accent5\"/></w14:shadow>' +\n  '<w14:textOutline w14:w=\"6731\" w14:cap=\"flat\" w14:cmpd=\"sng\" w14:algn=\"ctr\">' +\n  '<w14:solidFill><w14:schemeClr w14:val=\"bg1\"/></w14:solidFill>' +\n  '<w14:prstDash w14:val=\"solid\"/><w14:round/></w14:textOutline>';\n\nfunction run(text) {\n  const preserve = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : '';\n  return (\n    '<w:r><w:rPr>' + rPrInner + '</w:rPr>' +\n    '<w:t' + preserve + '>' + text + '</w:t></w:r>'\n  );\n}\n\nconst paragraphXml =\n  '<w:p>' +\n    '<w:pPr><w:jc w:val=\"center\"/><w:rPr>' + rPrInner + '</w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    run('C') +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    run('heck') +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    run(' ') +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    run('list') +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    run(' guia para entrevista') +\n  '</w:p>';\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document ' +\n          'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n          'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n          '<w:body>' + paragraphXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// Replace the whole body with the fully-formatted paragraph. insertOoxml\n// keeps the document's existing section properties (sectPr) intact, which\n// matches the Word JS API's behaviour (Body excludes section marks).\nconst body = context.document.body;\nbody.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# Source paragraph:      \"F\u00c9RIAS\"  (plain, no formatting)\n# Target paragraph:      \"Check list guia para entrevista\"\n#   - paragraph centered\n#   - bold + italic, dark-grey themed color, size 14pt (sz=28 half-points)\n#   - WordArt-style text effects (shadow + outline) carried as w14 run\n#     extensions that have no classic Word-OM property, so the run /\n#     paragraph XML is written directly via Range.InsertXML().\n#   - a `_GoBack` bookmark sits between \"C\" and \"heck\" (an artifact of\n#     Word leaving the last-edit-position bookmark behind) and the two\n#     words \"Check\"/\"list\" are wrapped in proofErr spell-check markers.\n\n$d = $word.ActiveDocument\n\n# The run formatting (identical on every run and on the paragraph mark)\n# reproduced from the target OOXML.\n$rPrInner = '<w:b/><w:i/><w:iCs/>' +\n  '<w:color w:val=\"262626\" w:themeColor=\"text1\" w:themeTint=\"D9\"/>' +\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>' +\n  '<w14:shadow w14:blurRad=\"0\" w14:dist=\"38100\" w14:dir=\"2700000\" ' +\n  'w14:sx=\"100000\" w14:sy=\"100000\" w14:kx=\"0\" w14:ky=\"0\" w14:algn=\"bl\">' +\n  '<w14:schemeClr w14:val=\"accent5\"/></w14:shadow>' +\n  '<w14:textOutline w14:w=\"6731\" w14:cap=\"flat\" w14:cmpd=\"sng\" w14:algn=\"ctr\">' +\n  '<w14:solidFill><w14:schemeClr w14:val=\"bg1\"/></w14:solidFill>' +\n  '<w14:prstDash w14:val=\"solid\"/><w14:round/></w14:textOutline>'\n\nfunction New-Run([string]$text) {\n    $preserve = ''\n    if ($text -match '^\\s' -or $text -match '\\s$') {\n        $preserve = ' xml:space=\"preserve\"'\n    }\n    return '<w:r><w:rPr>' + $rPrInner + '</w:rPr><w:t' + $preserve + '>' + $text + '</w:t></w:r>'\n}\n\n$paragraphXml =\n  '<w:p>' +\n    '<w:pPr><w:jc w:val=\"center\"/><w:rPr>' + $rPrInner + '</w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    (New-Run 'C') +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    (New-Run 'heck') +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    (New-Run ' ') +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    (New-Run 'list') +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    (New-Run ' guia para entrevista') +\n  '</w:p>'\n\n$ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document ' +\n          'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n          'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n          '<w:body>' + $paragraphXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n# Replace the whole document content with the fully-formatted paragraph.\n# InsertXML replaces the contents of the exact Range it is called on, and\n# keeps the document's existing section properties (sectPr) intact.\n$d.Content.InsertXML($ooxmlPackage)\n"}
